$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

# Read existing File Name / P_max pairs into an array of objects
$items = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fileName = $ws.Cells.Item($r, 1).Value()
    if ([string]::IsNullOrEmpty($fileName)) { continue }
    $pmax = $ws.Cells.Item($r, 2).Value()

    # Electrode location is the filename up to "_monopolar"
    $loc = $fileName.Split("_monopolar")[0]

    # Split the location into its leading letters and trailing number
    # so rows can be sorted alphabetically then numerically (A1, A3, ... A15, B13, ...)
    if ($loc -match '^([A-Za-z]+)(\d+)$') {
        $letters = $matches[1]
        $number = [int]$matches[2]
    } else {
        $letters = $loc
        $number = 0
    }

    $items += [PSCustomObject]@{
        FileName = $fileName
        Pmax     = $pmax
        Loc      = $loc
        Letters  = $letters
        Number   = $number
    }
}

# Sort rows by electrode location: letters first, then numeric part
$sorted = $items | Sort-Object Letters, Number

# Add the new "Electrode Locations" header, matching the style of the existing headers
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "Electrode Locations"

# Write the sorted data back, including the new Electrode Locations column
$row = 2
foreach ($item in $sorted) {
    $ws.Cells.Item($row, 1).Value = $item.FileName
    $ws.Cells.Item($row, 2).Value = $item.Pmax
    $ws.Cells.Item($row, 3).Value = $item.Loc
    $row++
}
